$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "20.240.39"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.69%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.441.97"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.38%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.54%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9191"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -8.25%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "274.84"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.59%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3635"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.06%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.12%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "38.71"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.66%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.019"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.22%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06483"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.00%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9985"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.34%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.318"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.07%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.41"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.05%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.026"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.18%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001007"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.28%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.441.52"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.77%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9354"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -6.65%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.05614"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.97%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "67.49"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -3.91%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.333"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -4.44%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.18"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -3.30%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.72"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.70%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.239"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.12%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "20.273.55"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.80%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "139.21"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.95%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.051"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -8.60%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.88"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.13%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.593.54"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.09%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "110.06"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.73%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.006"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.35%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.821"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -9.11%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7829"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.44%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07650"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.18%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.456"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.44%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05766"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.11%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +5.00%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.633"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.73%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01986"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.49%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.04%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1838"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.41%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9269"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -7.40%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.019"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -16.37%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5182"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.62%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.477"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.89%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "11.68"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -4.76%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "115.67"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.44%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5084"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.63%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.726"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.94%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06355"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.83%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9916"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.97%  "
